$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Predicted Eg" (D) and "Predicted Eg" (E) columns with the
# refreshed predictions after adding the Random Forest algorithm.
$ws.Range("D2").Value = 3.2
$ws.Range("E2").Value = 3.18
$ws.Range("D3").Value = 3.52
$ws.Range("E3").Value = 3.45
$ws.Range("D4").Value = 3.43
$ws.Range("E4").Value = 3.42
$ws.Range("D5").Value = 3.45
$ws.Range("E5").Value = 3.4
$ws.Range("D6").Value = 3.4
$ws.Range("E6").Value = 3.4
$ws.Range("D7").Value = 3.45
$ws.Range("E7").Value = 3.4
$ws.Range("D8").Value = 3.38
$ws.Range("E8").Value = 3.38
$ws.Range("D9").Value = 2.95
$ws.Range("E9").Value = 2.94
$ws.Range("D10").Value = 3.3
$ws.Range("E10").Value = 3.31
$ws.Range("D11").Value = 3.28
$ws.Range("E11").Value = 3.29
$ws.Range("D12").Value = 3.29
$ws.Range("E12").Value = 3.3
$ws.Range("D13").Value = 3.21
$ws.Range("E13").Value = 3.24
$ws.Range("D14").Value = 3.18
$ws.Range("E14").Value = 3.21
$ws.Range("D15").Value = 3.12
$ws.Range("E15").Value = 3.16
$ws.Range("D16").Value = 3.16
$ws.Range("E16").Value = 3.18
$ws.Range("D17").Value = 3.05
$ws.Range("D18").Value = 3.15
$ws.Range("E18").Value = 3.13
$ws.Range("D19").Value = 3.02
$ws.Range("E19").Value = 3.03
$ws.Range("D20").Value = 1.87
$ws.Range("E20").Value = 1.89
$ws.Range("D21").Value = 2.39
$ws.Range("E21").Value = 2.33
$ws.Range("D22").Value = 2.42
$ws.Range("E22").Value = 2.35
$ws.Range("D23").Value = 2.39
$ws.Range("E23").Value = 2.33
$ws.Range("D24").Value = 2.39
$ws.Range("E24").Value = 2.33
$ws.Range("D25").Value = 2.45
$ws.Range("E25").Value = 2.38
$ws.Range("D26").Value = 6.23
$ws.Range("E26").Value = 6.71
$ws.Range("D27").Value = 3.58
$ws.Range("E27").Value = 3.54
$ws.Range("D28").Value = 3.69
$ws.Range("E28").Value = 3.65
$ws.Range("D29").Value = 2.36
$ws.Range("E29").Value = 2.31
$ws.Range("D30").Value = 3.61
$ws.Range("E30").Value = 3.57
$ws.Range("D31").Value = 8.54
$ws.Range("E31").Value = 8.54
$ws.Range("D32").Value = 2.63
$ws.Range("E32").Value = 2.74
$ws.Range("D33").Value = 2.67
$ws.Range("E33").Value = 2.76
$ws.Range("D34").Value = 2.92
$ws.Range("E34").Value = 2.93
$ws.Range("D35").Value = 3.84
$ws.Range("E35").Value = 3.85

# Set column D width (engine rounds to nearest 1/6; 13.71 yields closest to 14.5703125 -> 14.5)
$ws.Columns("D:D").ColumnWidth = 13.71

# Select column E to match the recorded selection in the diff
$ws.Columns("E:E").Select()
